$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), pushing all existing
# price rows down by one, and copy formatting from the row below so the
# new row matches the rest of the table.
$ws.Rows.Item(2).Insert()
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row with the latest Hindalco circular data.
$ws.Cells.Item(2, 1).Value2 = 39
$ws.Cells.Item(2, 2).Value2 = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Cells.Item(2, 3).Value2 = "P1020"
$ws.Cells.Item(2, 4).Value2 = 270.25
$ws.Cells.Item(2, 5).Value2 = "30.09.2025"
$ws.Cells.Item(2, 6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-30-september-2025.pdf"

# Row 20 (previously row 19) now carries a circular link that needs its
# own hyperlink added, matching the other linked rows above it.
$ws.Hyperlinks.Add($ws.Cells.Item(20, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")
$ws.Range("F21").Copy()
$ws.Range("F20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "done"
